# Apply the "Add files via upload" edit to the Partidos sheet:
#  1. Rename a specific set of historical "Invitado" entries to "Montaño".
#  2. Append 11 new match rows (row 485-495) for the match played on 2025-10-18
#     (serial date 45948), one per player, across the Amarillo/Azul teams.
#  3. Leave the active cell selection on H498, matching the saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Partidos")

# --- 1. Rename "Invitado" -> "Montaño" for the specific rows that were edited ---
$renamedRows = @(363, 374, 382, 445, 456, 472)
foreach ($r in $renamedRows) {
    $ws.Cells.Item($r, 2).Value = "Montaño"
}

# --- 2. Append the new rows for the 2025-10-18 match ---
# Columns: A=fecha B=jugador C=equipo D=posicion E=goles F=autogoles
#          G=arquero H=goles_recibidos I=tarjetas_amarillas J=tarjetas_rojas
#          K=asistencias L=Penales_Atajados
$newRows = @(
    @(45948, "Alexander Uribe",          "Amarillo", "Mediocampista", 2, 0, $false, 0, 0, 0, 0, 0),
    @(45948, "Andres Tangarife",         "Amarillo", "Delantero",     1, 0, $false, 0, 0, 0, 0, 0),
    @(45948, "Jefferson Delgado",        "Amarillo", "Mediocampista", 1, 0, $false, 0, 0, 0, 2, 0),
    @(45948, "Francisco Javier Duran",   "Amarillo", "Defensa",       0, 0, $false, 0, 0, 0, 1, 0),
    @(45948, "Montaño",                  "Amarillo", "Arquero",       0, 0, $true,  6, 0, 0, 0, 0),
    @(45948, "Armando Murillo",          "Azul",     "Defensa",       2, 0, $false, 0, 0, 0, 0, 0),
    @(45948, "Carlos Fernando Valencia", "Azul",     "Delantero",     2, 0, $false, 0, 0, 0, 1, 0),
    @(45948, "Juan David Espinal",       "Azul",     "Mediocampista", 2, 0, $false, 0, 0, 0, 0, 0),
    @(45948, "Bryan Andres Burgos",      "Azul",     "Mediocampista", 0, 0, $false, 0, 0, 0, 3, 0),
    @(45948, "Edwin Casas",              "Azul",     "Mediocampista", 0, 0, $false, 0, 0, 0, 1, 0),
    @(45948, "Gember Marin Sarria",      "Azul",     "Arquero",       0, 0, $true,  4, 0, 0, 0, 0)
)

$startRow = 485
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($row, 1).Value = $data[0]
    $ws.Cells.Item($row, 2).Value = $data[1]
    $ws.Cells.Item($row, 3).Value = $data[2]
    $ws.Cells.Item($row, 4).Value = $data[3]
    $ws.Cells.Item($row, 5).Value = $data[4]
    $ws.Cells.Item($row, 6).Value = $data[5]
    $ws.Cells.Item($row, 7).Value = $data[6]
    $ws.Cells.Item($row, 8).Value = $data[7]
    $ws.Cells.Item($row, 9).Value = $data[8]
    $ws.Cells.Item($row, 10).Value = $data[9]
    $ws.Cells.Item($row, 11).Value = $data[10]
    $ws.Cells.Item($row, 12).Value = $data[11]
}

# --- 3. Restore the active selection recorded in the saved view ---
$ws.Range("H498").Select()
